$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a space value in A19 (matching the styled look of the rest of column A)
# and in B20 (matching the unstyled look of the rest of column B).
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = " "

$ws.Range("B20").Value = " "

# Update the active selection to A19, matching the saved view state.
$ws.Range("A19").Select()
